# Update "想去人数" (column F) counts for several con/expo events on the
# "展览" (sheet1) and "全部类型" (sheet4) sheets, per output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# Map: sheet name -> list of (row, newValue)
$updates = @{
    "展览"     = @{ 2 = 68; 3 = 800; 6 = 102; 8 = 4106; 10 = 4824; 11 = 535; 12 = 1209 }
    "全部类型" = @{ 2 = 68; 3 = 800; 6 = 102; 9 = 4106; 11 = 4824; 12 = 535; 13 = 1209 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $newVal = $rowsMap[$row]
        $ws.Range("F$row").Value = $newVal
    }
}
